$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -38.461538461538
$ws.Range("I16").Value = 89
$ws.Range("J16").Value = 91
$ws.Range("K16").Value = -2.197802197802
$ws.Range("L16").Value = 28.985507246376
$ws.Range("M16").Value = 97.777777777777
$ws.Range("N16").Value = -83.333333333333
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 27.272727272727
$ws.Range("I17").Value = 90
$ws.Range("J17").Value = 87
$ws.Range("K17").Value = 3.448275862068
$ws.Range("L17").Value = 21.621621621621
$ws.Range("M17").Value = 95.652173913043
$ws.Range("N17").Value = -23.728813559322
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 174
$ws.Range("K18").Value = -32.183908045977
$ws.Range("L18").Value = 38.823529411764
$ws.Range("M18").Value = 5.357142857142
$ws.Range("N18").Value = -78.228782287822
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -30.434782608695
$ws.Range("F19").Value = 104
$ws.Range("G19").Value = 104
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 764
$ws.Range("J19").Value = 780
$ws.Range("K19").Value = -2.051282051282
$ws.Range("L19").Value = 69.026548672566
$ws.Range("M19").Value = 13.353115727003
$ws.Range("N19").Value = -69.11883589329
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 7
$ws.Range("I20").Value = 44
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 37.5
$ws.Range("M20").Value = 76
$ws.Range("N20").Value = -92.321116928446
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -16.129032258064
$ws.Range("F21").Value = 145
$ws.Range("G21").Value = 152
$ws.Range("H21").Value = -4.605263157894
$ws.Range("I21").Value = 1115
$ws.Range("J21").Value = 1190
$ws.Range("K21").Value = -6.302521008403
$ws.Range("L21").Value = 53.793103448275
$ws.Range("M21").Value = 23.204419889502
$ws.Range("N21").Value = -73.813997181775
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 54
$ws.Range("J22").Value = 67
$ws.Range("K22").Value = -19.402985074626
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = 31.70731707317
$ws.Range("C24").Value = 91
$ws.Range("D24").Value = 73
$ws.Range("E24").Value = 24.657534246575
$ws.Range("F24").Value = 358
$ws.Range("G24").Value = 297
$ws.Range("H24").Value = 20.53872053872
$ws.Range("I24").Value = 2502
$ws.Range("J24").Value = 2497
$ws.Range("K24").Value = 0.200240288346
$ws.Range("L24").Value = 101.611603545528
$ws.Range("M24").Value = 125.812274368231
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 85.714285714285
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 68.965517241379
$ws.Range("I25").Value = 253
$ws.Range("J25").Value = 226
$ws.Range("K25").Value = 11.946902654867
$ws.Range("L25").Value = 25.247524752475
$ws.Range("M25").Value = 59.119496855345
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = -43.75
$ws.Range("I27").Value = 61
$ws.Range("J27").Value = 74
$ws.Range("K27").Value = -17.567567567567
$ws.Range("L27").Value = 22

# --- Row 26: D26/E26 switch from numbers to text (shared strings "0" / "***.*") ---
# Force text entry with a leading apostrophe so COM does not re-parse as numbers,
# then paste the number-format/style from a neighboring text cell (C26) so the
# resulting style index matches the other "n/a" cells in the row.
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "'***.*"
$ws.Range("C26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").PasteSpecial(-4122)

# --- Row 27: D27/E27 switch from text back to numbers ---
# Set the numeric values, then paste formats from neighboring numeric cells
# (F27 -> integer style, H27 -> percent-change style) to match target styles.
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$excel.CutCopyMode = 0

